$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the index bug: add the "Initial_value" reference column with
# the values from the template.
$ws.Range("D1").Value = "Initial_value"
$ws.Range("D2").Value = 4278.04
$ws.Range("D3").Value = 3170.73
$ws.Range("D4").Value = 1268.41

$ws.Range("D4").Select()
